$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 values - B2 becomes blank, C2/D2/E2 updated
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 1.2093452929750623
$ws.Range("D2").Value = 0.27049389530226514
$ws.Range("E2").Value = 3.4592388228240232

# Row 3 values
$ws.Range("B3").Value = 0.63775836925333151
$ws.Range("C3").Value = 1.8997263969175724
$ws.Range("D3").Value = 0.40143632450453087
$ws.Range("E3").Value = 1.6867536660769435

# Update selection to match new data extents
$ws.Range("B1:E3").Select()
